$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=8; Tag='ba'; Desc='Appreciation'}
    @{Row=13; Tag='sv'; Desc='Statement-opinion'}
    @{Row=15; Tag='sv'; Desc='Statement-opinion'}
    @{Row=20; Tag='sv'; Desc='Statement-opinion'}
    @{Row=35; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=46; Tag='%'; Desc='Uninterpretable'}
    @{Row=51; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=62; Tag='sv'; Desc='Statement-opinion'}
    @{Row=64; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=67; Tag='ba'; Desc='Appreciation'}
    @{Row=71; Tag='sv'; Desc='Statement-opinion'}
    @{Row=73; Tag='sv'; Desc='Statement-opinion'}
    @{Row=83; Tag='aa'; Desc='Agree/Accept'}
    @{Row=95; Tag='aa'; Desc='Agree/Accept'}
    @{Row=96; Tag='aa'; Desc='Agree/Accept'}
    @{Row=102; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=106; Tag='%'; Desc='Uninterpretable'}
    @{Row=111; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=124; Tag='aa'; Desc='Agree/Accept'}
    @{Row=130; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=139; Tag='sv'; Desc='Statement-opinion'}
    @{Row=156; Tag='aa'; Desc='Agree/Accept'}
    @{Row=159; Tag='aa'; Desc='Agree/Accept'}
    @{Row=162; Tag='aa'; Desc='Agree/Accept'}
    @{Row=169; Tag='aa'; Desc='Agree/Accept'}
    @{Row=173; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=188; Tag='sv'; Desc='Statement-opinion'}
    @{Row=189; Tag='aa'; Desc='Agree/Accept'}
    @{Row=195; Tag='sv'; Desc='Statement-opinion'}
    @{Row=200; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=214; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=229; Tag='sv'; Desc='Statement-opinion'}
    @{Row=256; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=274; Tag='sv'; Desc='Statement-opinion'}
    @{Row=276; Tag='qy'; Desc='Yes-No-Question'}
    @{Row=277; Tag='sv'; Desc='Statement-opinion'}
    @{Row=282; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=283; Tag='aa'; Desc='Agree/Accept'}
    @{Row=293; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=302; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=316; Tag='sv'; Desc='Statement-opinion'}
    @{Row=330; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=331; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=341; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=344; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=367; Tag='sv'; Desc='Statement-opinion'}
    @{Row=372; Tag='ba'; Desc='Appreciation'}
    @{Row=391; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=405; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=422; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=425; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=429; Tag='sv'; Desc='Statement-opinion'}
    @{Row=436; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=464; Tag='ba'; Desc='Appreciation'}
    @{Row=481; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=485; Tag='sv'; Desc='Statement-opinion'}
    @{Row=503; Tag='aa'; Desc='Agree/Accept'}
    @{Row=504; Tag='aa'; Desc='Agree/Accept'}
    @{Row=508; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=511; Tag='aa'; Desc='Agree/Accept'}
    @{Row=516; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=522; Tag='ba'; Desc='Appreciation'}
    @{Row=540; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=547; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=549; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=579; Tag='sv'; Desc='Statement-opinion'}
    @{Row=582; Tag='sv'; Desc='Statement-opinion'}
    @{Row=586; Tag='ba'; Desc='Appreciation'}
    @{Row=605; Tag='sv'; Desc='Statement-opinion'}
    @{Row=627; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=632; Tag='ba'; Desc='Appreciation'}
    @{Row=634; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=664; Tag='aa'; Desc='Agree/Accept'}
    @{Row=676; Tag='sv'; Desc='Statement-opinion'}
    @{Row=686; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=689; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=695; Tag='qy'; Desc='Yes-No-Question'}
    @{Row=700; Tag='qy'; Desc='Yes-No-Question'}
    @{Row=706; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=711; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=729; Tag='aa'; Desc='Agree/Accept'}
    @{Row=732; Tag='sv'; Desc='Statement-opinion'}
    @{Row=733; Tag='sv'; Desc='Statement-opinion'}
    @{Row=763; Tag='sv'; Desc='Statement-opinion'}
    @{Row=783; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=811; Tag='aa'; Desc='Agree/Accept'}
    @{Row=825; Tag='ba'; Desc='Appreciation'}
    @{Row=871; Tag='aa'; Desc='Agree/Accept'}
    @{Row=886; Tag='aa'; Desc='Agree/Accept'}
    @{Row=889; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=890; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=892; Tag='aa'; Desc='Agree/Accept'}
    @{Row=893; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=909; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=913; Tag='sd'; Desc='Statement-non-opinion'}
    @{Row=931; Tag='ba'; Desc='Appreciation'}
    @{Row=933; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=934; Tag='%'; Desc='Uninterpretable'}
    @{Row=935; Tag='b'; Desc='Acknowledge (Backchannel)'}
    @{Row=936; Tag='aa'; Desc='Agree/Accept'}
    @{Row=939; Tag='aa'; Desc='Agree/Accept'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Desc
}
